$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '''69.734.49'
$ws.Cells.Item(2, 5).Value = '  -2.23%  '

$ws.Cells.Item(3, 4).Value = '''3.668.73'
$ws.Cells.Item(3, 5).Value = '  -2.86%  '

$ws.Cells.Item(4, 5).Value = '  +0.17%  '

$ws.Cells.Item(5, 4).Value = '''614.82'
$ws.Cells.Item(5, 5).Value = '  +0.05%  '

$ws.Cells.Item(6, 4).Value = '''177.16'
$ws.Cells.Item(6, 5).Value = '  -1.58%  '

$ws.Cells.Item(7, 4).Value = '''3.668.20'
$ws.Cells.Item(7, 5).Value = '  -2.74%  '

$ws.Cells.Item(8, 5).Value = '  +0.03%  '

$ws.Cells.Item(9, 4).Value = '''0.529'
$ws.Cells.Item(9, 5).Value = '  -2.58%  '

$ws.Cells.Item(10, 4).Value = '''0.163'
$ws.Cells.Item(10, 5).Value = '  -2.53%  '

$ws.Cells.Item(11, 4).Value = '''6.23'
$ws.Cells.Item(11, 5).Value = '  -3.59%  '

$ws.Cells.Item(12, 4).Value = '''0.478'
$ws.Cells.Item(12, 5).Value = '  -4.87%  '

$ws.Cells.Item(13, 4).Value = '''39.55'
$ws.Cells.Item(13, 5).Value = '  -2.61%  '

$ws.Cells.Item(14, 4).Value = '''0.0000252'
$ws.Cells.Item(14, 5).Value = '  -2.51%  '

$ws.Cells.Item(15, 4).Value = '''4.291.17'
$ws.Cells.Item(15, 5).Value = '  -2.61%  '

$ws.Cells.Item(16, 4).Value = '''3.685.47'
$ws.Cells.Item(16, 5).Value = '  -2.46%  '

$ws.Cells.Item(17, 4).Value = '''69.699.91'
$ws.Cells.Item(17, 5).Value = '  -2.39%  '

$ws.Cells.Item(18, 4).Value = '''0.121'
$ws.Cells.Item(18, 5).Value = '  -2.12%  '

$ws.Cells.Item(19, 4).Value = '''7.47'
$ws.Cells.Item(19, 5).Value = '  -0.97%  '

$ws.Cells.Item(20, 4).Value = '''16.27'
$ws.Cells.Item(20, 5).Value = '  -2.77%  '

$ws.Cells.Item(21, 4).Value = '''498.76'
$ws.Cells.Item(21, 5).Value = '  -4.91%  '

$ws.Cells.Item(22, 4).Value = '''9.11'
$ws.Cells.Item(22, 5).Value = '  -2.51%  '

$ws.Cells.Item(23, 4).Value = '''0.707'
$ws.Cells.Item(23, 5).Value = '  -5.46%  '

$ws.Cells.Item(24, 4).Value = '''2.55'
$ws.Cells.Item(24, 5).Value = '  +1.70%  '

$ws.Cells.Item(25, 4).Value = '''85.27'
$ws.Cells.Item(25, 5).Value = '  -3.67%  '

$ws.Cells.Item(26, 4).Value = '''11.30'
$ws.Cells.Item(26, 5).Value = '  +2.12%  '

$ws.Cells.Item(27, 4).Value = '''12.92'
$ws.Cells.Item(27, 5).Value = '  -4.63%  '

$ws.Cells.Item(28, 4).Value = '''0.0000130'
$ws.Cells.Item(28, 5).Value = '  +8.66%  '

$ws.Cells.Item(29, 5).Value = '  -0.28%  '

$ws.Cells.Item(30, 2).Value = 'PancakeSwap'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Cells.Item(30, 4).Value = '''2.90'
$ws.Cells.Item(30, 5).Value = '  -0.61%  '

$ws.Cells.Item(31, 2).Value = 'ImmutableX'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Cells.Item(31, 4).Value = '''2.43'
$ws.Cells.Item(31, 5).Value = '  -3.82%  '

$ws.Cells.Item(32, 4).Value = '''7.78'
$ws.Cells.Item(32, 5).Value = '  -3.43%  '

$ws.Cells.Item(33, 4).Value = '''29.87'
$ws.Cells.Item(33, 5).Value = '  -7.04%  '

$ws.Cells.Item(34, 4).Value = '''0.113'
$ws.Cells.Item(34, 5).Value = '  -2.14%  '

$ws.Cells.Item(35, 5).Value = '  +0.08%  '

$ws.Cells.Item(36, 4).Value = '''1.04'
$ws.Cells.Item(36, 5).Value = '  -2.03%  '

$ws.Cells.Item(37, 4).Value = '''6.00'
$ws.Cells.Item(37, 5).Value = '  -2.30%  '

$ws.Cells.Item(38, 4).Value = '''0.137'
$ws.Cells.Item(38, 5).Value = '  +2.97%  '

$ws.Cells.Item(39, 4).Value = '''0.334'
$ws.Cells.Item(39, 5).Value = '  -2.27%  '

$ws.Cells.Item(40, 4).Value = '''2.05'
$ws.Cells.Item(40, 5).Value = '  -8.65%  '

$ws.Cells.Item(41, 2).Value = 'OKB'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Cells.Item(41, 4).Value = '''49.73'
$ws.Cells.Item(41, 5).Value = '  -4.32%  '

$ws.Cells.Item(42, 2).Value = 'Arweave'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Cells.Item(42, 4).Value = '''45.14'
$ws.Cells.Item(42, 5).Value = '  +2.74%  '

$ws.Cells.Item(43, 4).Value = '''423.62'
$ws.Cells.Item(43, 5).Value = '  -1.85%  '

$ws.Cells.Item(44, 4).Value = '''2.87'
$ws.Cells.Item(44, 5).Value = '  +2.30%  '

$ws.Cells.Item(45, 4).Value = '''8.50'
$ws.Cells.Item(45, 5).Value = '  -3.68%  '

$ws.Cells.Item(46, 4).Value = '''2.935.58'
$ws.Cells.Item(46, 5).Value = '  -7.03%  '

$ws.Cells.Item(47, 4).Value = '''0.0356'
$ws.Cells.Item(47, 5).Value = '  -3.17%  '

$ws.Cells.Item(48, 2).Value = 'USDe'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Cells.Item(48, 4).Value = '''1.00'
$ws.Cells.Item(48, 5).Value = '  +0.00%  '

$ws.Cells.Item(49, 2).Value = 'InjectiveProtocol'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Cells.Item(49, 4).Value = '''27.02'
$ws.Cells.Item(49, 5).Value = '  -3.33%  '

$ws.Cells.Item(50, 4).Value = '''135.80'
$ws.Cells.Item(50, 5).Value = '  -3.56%  '

$ws.Cells.Item(51, 4).Value = '''2.43'
$ws.Cells.Item(51, 5).Value = '  -1.25%  '
